$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 31 contain a date value (stored as serial 46075)
# that should be incremented by one day to 46076.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 46075) {
        $cell.Value2 = 46076
    }
}
